$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['HardwareFault', 'SoftwareFault']"

$ws.Range("D27").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['HardwareFault', 'SoftwareFault']"

$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal']"
